# Apply content edits to Sheet1 of the workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Column A: wording / capitalization fixes -------------------------------
$ws.Range("A2").Value  = "Inflation and real returns"
$ws.Range("A3").Value  = "Net present value"
$ws.Range("A4").Value  = "Internal rate of return"
$ws.Range("A6").Value  = "Amortization schedule"
$ws.Range("A8").Value  = "Bond prices and yields"
$ws.Range("A11").Value = "Treasury inflation protected securities"
$ws.Range("A56").Value = "Option portfolios"

# --- Column B: hyperlink display text / target, github.com -> githubtocolab.com ---
$linkRows = @(1,2,3,4,5,9,11,16,17,23,24,27,29,30,37,38,39,40,42,45,56,60,61,75)
foreach ($r in $linkRows) {
    $cell = $ws.Range("B$r")
    $oldUrl = $cell.Value2
    $newUrl = $oldUrl -replace "https://github\.com/", "https://githubtocolab.com/"
    $cell.Value = $newUrl
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address($false, $false) -eq "B$r") {
            $hl.Address = $newUrl
            $hl.TextToDisplay = $newUrl
        }
    }
}

# --- Sheet view: move the active selection to A40 ---------------------------
$ws.Activate()
$ws.Range("A40").Select()
